$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "ECs" sending-cluster rows (original rows 2-4); the
# remaining FAPs/MuSCs sending rows shift up and become the new rows 2-7.
$ws.Rows("2:4").Delete()

# Updated (new TPM) numeric values for columns E..T on the new rows 2-7.
# Row 2: FAPs -> Lif/Il6st -> ECs
# Row 3: FAPs -> Lif/Il6st -> FAPs
# Row 4: FAPs -> Lif/Il6st -> MuSCs
# Row 5: MuSCs -> Lif/Il6st -> ECs
# Row 6: MuSCs -> Lif/Il6st -> FAPs
# Row 7: MuSCs -> Lif/Il6st -> MuSCs
$newValues = @{
    2 = @(3, 1, 2.335066666666667, 7.0052, 0.9647765390673073, 0.9647765390673074, 3, 1, 19.827687, 59.483061, 0.1538389073329896, 0.1538389073329896, 46.2989709908, 416.6907389172, 0.1484201685906179, 0.1484201685906179)
    3 = @(3, 1, 2.335066666666667, 7.0052, 0.9647765390673073, 0.9647765390673074, 3, 1, 85.11961100000001, 255.358833, 0.6604253914664442, 0.6604253914664441, 198.7599663257333, 1788.8396969316, 0.6371629234911677, 0.6371629234911675)
    4 = @(3, 1, 2.335066666666667, 7.0052, 0.9647765390673073, 0.9647765390673074, 3, 1, 23.93873833333333, 71.816215, 0.1857357012005663, 0.1857357012005663, 55.89854992422222, 503.086949318, 0.1791934469855219, 0.1791934469855219)
    5 = @(1, 0.3333333333333333, 0.08525199999999999, 0.255756, 0.03522346093269261, 0.03522346093269261, 3, 1, 19.827687, 59.483061, 0.1538389073329896, 0.1538389073329896, 1.690349972124, 15.213149749116, 0.005418738742371679, 0.005418738742371678)
    6 = @(1, 0.3333333333333333, 0.08525199999999999, 0.255756, 0.03522346093269261, 0.03522346093269261, 3, 1, 85.11961100000001, 255.358833, 0.6604253914664442, 0.6604253914664441, 7.256617076972, 65.309553692748, 0.02326246797527652, 0.02326246797527651)
    7 = @(1, 0.3333333333333333, 0.08525199999999999, 0.255756, 0.03522346093269261, 0.03522346093269261, 3, 1, 23.93873833333333, 71.816215, 0.1857357012005663, 0.1857357012005663, 2.040825320393333, 18.36742788354, 0.006542254215044415, 0.006542254215044415)
}

foreach ($r in $newValues.Keys) {
    $vals = $newValues[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 5).Value = $vals[$i]
    }
}
